$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 238, shifting the existing rows 238-286 down to 239-287
$ws.Rows(238).Insert()

# Populate the newly inserted row 238 with the new weekly price record
$ws.Cells.Item(238,1).Value  = 11
$ws.Cells.Item(238,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(238,3).Value  = "Bíobío"
$ws.Cells.Item(238,4).Value  = 44855
$ws.Cells.Item(238,5).Value  = 8
$ws.Cells.Item(238,6).Value  = 100114013
$ws.Cells.Item(238,7).Value  = "Zanahoria"
$ws.Cells.Item(238,8).Value  = "Sin especificar"
$ws.Cells.Item(238,9).Value  = "Primera"
$ws.Cells.Item(238,10).Value = 1000
$ws.Cells.Item(238,11).Value = 21000
$ws.Cells.Item(238,12).Value = 22000
$ws.Cells.Item(238,13).Value = 21500
$ws.Cells.Item(238,14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(238,15).Value = "Región de La Araucanía"
$ws.Cells.Item(238,16).Value = 1075
$ws.Cells.Item(238,17).Value = 20
$ws.Cells.Item(238,18).Value = "Hortaliza"
